# Re-ran resolve and classify+summarise steps after changes to mapping file.
# For Chandigarh, the refreshed mapping produced no matching species, so the
# "Range Status" breakdown collapses to zero counts (and the now-undefined
# percentage column is dropped), the "Range Analysis" species count in
# "Species qualification" drops to 0, while the "High Priority break-up"
# sheet gains its "New High Species" figures.

$wb = $excel.ActiveWorkbook

# --- Sheet: Range Status ---
$wsRange = $wb.Worksheets.Item("Range Status")

$wsRange.Range("B2").Value = 0
$wsRange.Range("C2").ClearContents()

$wsRange.Range("C3").ClearContents()

$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()

$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()

$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()

$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# --- Sheet: Species qualification ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# --- Sheet: High Priority break-up ---
$wsBreak = $wb.Worksheets.Item("High Priority break-up")
$wsBreak.Range("D2").Value = 3
$wsBreak.Range("E2").Value = 100
